# Executing new implementation of readFile().
#
# The original run of the tool populated RQ1-RQ4 with data read from the
# repositories. This re-run adds four new sheets (RQ1a-RQ4a) holding the
# results of a new readFile() implementation; some of the source files
# could not be (fully) re-read, so several input cells in the new sheets
# are blank (the dependent ratio/percentage formulas consequently show
# #DIV/0!).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Selection / active-cell bookkeeping on the four pre-existing sheets.
# ---------------------------------------------------------------------
$rq1 = $wb.Worksheets.Item("RQ1")
$rq1.Range("C3").Select()

$rq2 = $wb.Worksheets.Item("RQ2")
$rq2.Range("E3").Select()

$rq3 = $wb.Worksheets.Item("RQ3")
$rq3.Range("A1:G8").Select()

$rq4 = $wb.Worksheets.Item("RQ4")
$rq4.Range("A1:G8").Select()

# ---------------------------------------------------------------------
# 2. RQ1a - copy of RQ1 with most of the per-project inputs cleared;
#    only the "Revisions" counts for Ctags and Carol survived the re-run.
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$rq1.Copy($null, $last)
$rq1a = $wb.Worksheets.Item($wb.Worksheets.Count)
$rq1a.Name = "RQ1a"

$rq1a.Range("C3:D3").ClearContents()
$rq1a.Range("F3:G3").ClearContents()
$rq1a.Range("C4:D4").ClearContents()
$rq1a.Range("F4:G4").ClearContents()
$rq1a.Range("C6:D6").ClearContents()
$rq1a.Range("F6:G6").ClearContents()
$rq1a.Range("C7:D7").ClearContents()
$rq1a.Range("F7:G7").ClearContents()

$rq1a.Range("G3:G7").Select()

# ---------------------------------------------------------------------
# 3. RQ2a - copy of RQ2; only the Ctags ("Brlcad" row untouched) row of
#    counts came back from the new readFile(), with updated numbers.
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$rq2.Copy($null, $last)
$rq2a = $wb.Worksheets.Item($wb.Worksheets.Count)
$rq2a.Name = "RQ2a"

$rq2a.Range("B3").Value = 82
$rq2a.Range("C3").Value = 88
$rq2a.Range("E3").Value = 10
$rq2a.Range("F3").Value = 10

$rq2a.Range("B4:C4").ClearContents()
$rq2a.Range("E4:F4").ClearContents()
$rq2a.Range("B6:C6").ClearContents()
$rq2a.Range("E6:F6").ClearContents()
$rq2a.Range("B7:C7").ClearContents()
$rq2a.Range("E7:F7").ClearContents()

$rq2a.Range("C4").Select()

# ---------------------------------------------------------------------
# 4. RQ3a - copy of RQ3 with every per-project input blanked out.
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$rq3.Copy($null, $last)
$rq3a = $wb.Worksheets.Item($wb.Worksheets.Count)
$rq3a.Name = "RQ3a"

$rq3a.Range("B3:C3").ClearContents()
$rq3a.Range("E3:F3").ClearContents()
$rq3a.Range("B4:C4").ClearContents()
$rq3a.Range("E4:F4").ClearContents()
$rq3a.Range("B6:C6").ClearContents()
$rq3a.Range("E6:F6").ClearContents()
$rq3a.Range("B7:C7").ClearContents()
$rq3a.Range("E7:F7").ClearContents()

$rq3a.Range("B3").Select()

# ---------------------------------------------------------------------
# 5. RQ4a - copy of RQ4 with every per-project input blanked out.
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$rq4.Copy($null, $last)
$rq4a = $wb.Worksheets.Item($wb.Worksheets.Count)
$rq4a.Name = "RQ4a"

$rq4a.Range("B3:C3").ClearContents()
$rq4a.Range("E3:F3").ClearContents()
$rq4a.Range("B4:C4").ClearContents()
$rq4a.Range("E4:F4").ClearContents()
$rq4a.Range("B6:C6").ClearContents()
$rq4a.Range("E6:F6").ClearContents()
$rq4a.Range("B7:C7").ClearContents()
$rq4a.Range("E7:F7").ClearContents()

$rq4a.Range("E4").Select()

# ---------------------------------------------------------------------
# 6. RQ2a ends up as the active tab (matches the saved window state).
# ---------------------------------------------------------------------
$rq2a.Activate()
